# Weekly data update for "Vega Modelo de Temuco - Zanahoria"
# Two new price records are inserted at rows 542-543 (pushing the existing
# rows 542-554 down to 544-556), matching a new weekly extract being
# prepended to the top of this price block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 542, shifting existing rows 542:554 down to 544:556
$ws.Rows("542:543").Insert()

# --- New row 542 ---
$ws.Cells.Item(542,1).Value = 10
$ws.Cells.Item(542,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(542,3).Value = "La Araucanía"
$ws.Cells.Item(542,4).Value = 45239
$ws.Cells.Item(542,5).Value = 9
$ws.Cells.Item(542,6).Value = 100114013
$ws.Cells.Item(542,7).Value = "Zanahoria"
$ws.Cells.Item(542,8).Value = "Sin especificar"
$ws.Cells.Item(542,9).Value = "Primera"
$ws.Cells.Item(542,10).Value = 150
$ws.Cells.Item(542,11).Value = 7000
$ws.Cells.Item(542,12).Value = 7000
$ws.Cells.Item(542,13).Value = 7000
$ws.Cells.Item(542,14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(542,15).Value = "Provincia del Elquí"
$ws.Cells.Item(542,16).Value = 350
$ws.Cells.Item(542,17).Value = 20
$ws.Cells.Item(542,18).Value = "Hortaliza"

# --- New row 543 ---
$ws.Cells.Item(543,1).Value = 10
$ws.Cells.Item(543,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(543,3).Value = "La Araucanía"
$ws.Cells.Item(543,4).Value = 45239
$ws.Cells.Item(543,5).Value = 9
$ws.Cells.Item(543,6).Value = 100114013
$ws.Cells.Item(543,7).Value = "Zanahoria"
$ws.Cells.Item(543,8).Value = "Sin especificar"
$ws.Cells.Item(543,9).Value = "Primera"
$ws.Cells.Item(543,10).Value = 200
$ws.Cells.Item(543,11).Value = 7000
$ws.Cells.Item(543,12).Value = 7000
$ws.Cells.Item(543,13).Value = 7000
$ws.Cells.Item(543,14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(543,15).Value = "Región del Maule"
$ws.Cells.Item(543,16).Value = 350
$ws.Cells.Item(543,17).Value = 20
$ws.Cells.Item(543,18).Value = "Hortaliza"
